$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table "Tabela1" currently covers A1:J70 (header + 69 data rows).
# Add one more data row to the table so the table range, autofilter and
# worksheet dimension all grow to include row 71 (matching native Excel
# "Insert table row" behaviour).
$lo = $ws.ListObjects.Item(1)
$newListRow = $lo.ListRows.Add()
$r = $newListRow.Range.Row

# Copy the formatting of an existing, unshaded data row (row 20 - column A
# uses the custom date format, column B uses a thousands-separated number
# format, and columns C:J use the general number format, all without the
# banding fill/border used on the even data rows) down onto the freshly
# inserted row so the new cells pick up the same formatting used
# throughout the table body instead of inheriting the column's default
# (text) format.
$ws.Range("A20:J20").Copy()
$ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 10)).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data for 20/5/2020 (Excel serial date 43971).
$ws.Cells.Item($r, 1).Value = 43971
$ws.Cells.Item($r, 2).Value = 72860
$ws.Cells.Item($r, 3).Value = 909
$ws.Cells.Item($r, 4).Value = 1468
$ws.Cells.Item($r, 5).Value = 0
$ws.Cells.Item($r, 6).Value = 21
$ws.Cells.Item($r, 7).Value = 3
$ws.Cells.Item($r, 8).Value = 2
$ws.Cells.Item($r, 9).Value = 106
$ws.Cells.Item($r, 10).Value = 1

# Match the saved selection state (the newly added row is selected).
$newRowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 10))
$newRowRange.Select()
